$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the E and F column header labels (C_A <-> A_C)
$ws.Range("E1").Value = "A_C"
$ws.Range("F1").Value = "C_A"

# Swap the E2 and F2 data values to match
$ws.Range("E2").Value = 8.838520634608827
$ws.Range("F2").Value = 0.0884978075470728
